# Scheduled runner update: refresh cached Universalis market-price
# snapshots (currentAveragePrice* / LevePrice* / LeveProfit* columns)
# across the per-job Leve profitability sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW). Values below mirror the latest market data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 10870482
$ws.Range("I112").Value = 796.6667
$ws.Range("J112").Value = 11628833
$ws.Range("K112").Value = 2390.0001
$ws.Range("L112").Value = 34886499
$ws.Range("M112").Value = -1282.0001
$ws.Range("N112").Value = -34888715
$ws.Range("H127").Value = 1827.1428
$ws.Range("I127").Value = 735.6
$ws.Range("J127").Value = 2168.25
$ws.Range("K127").Value = 2206.8
$ws.Range("L127").Value = 6504.75
$ws.Range("M127").Value = 2753.2
$ws.Range("N127").Value = -16424.75
$ws.Range("H132").Value = 2080.111
$ws.Range("I132").Value = 1362.9667
$ws.Range("J132").Value = 5665.8335
$ws.Range("K132").Value = 4088.9001
$ws.Range("L132").Value = 16997.5005
$ws.Range("M132").Value = -1558.9001
$ws.Range("N132").Value = -22057.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3454.4458
$ws.Range("I32").Value = 3010.1538
$ws.Range("J32").Value = 5058.8335
$ws.Range("K32").Value = 3010.1538
$ws.Range("L32").Value = 5058.8335
$ws.Range("M32").Value = -2723.1538
$ws.Range("N32").Value = -5632.8335
$ws.Range("H61").Value = 2363.8235
$ws.Range("I61").Value = 2145.04
$ws.Range("J61").Value = 2971.5557
$ws.Range("K61").Value = 2145.04
$ws.Range("L61").Value = 2971.5557
$ws.Range("M61").Value = -1933.04
$ws.Range("N61").Value = -3395.5557
$ws.Range("H74").Value = 1681.4595
$ws.Range("I74").Value = 1358.5834
$ws.Range("J74").Value = 2277.5386
$ws.Range("K74").Value = 1358.5834
$ws.Range("L74").Value = 2277.5386
$ws.Range("M74").Value = -484.5834
$ws.Range("N74").Value = -4025.5386
$ws.Range("H77").Value = 1681.4595
$ws.Range("I77").Value = 1358.5834
$ws.Range("J77").Value = 2277.5386
$ws.Range("K77").Value = 6792.916999999999
$ws.Range("L77").Value = 11387.693
$ws.Range("M77").Value = -2424.916999999999
$ws.Range("N77").Value = -20123.693
$ws.Range("H132").Value = 2231.9473
$ws.Range("I132").Value = 2071.3547
$ws.Range("J132").Value = 2943.1428
$ws.Range("K132").Value = 6214.0641
$ws.Range("L132").Value = 8829.428400000001
$ws.Range("M132").Value = -3684.0641
$ws.Range("N132").Value = -13889.4284
$ws.Range("H136").Value = 2363.8235
$ws.Range("I136").Value = 2145.04
$ws.Range("J136").Value = 2971.5557
$ws.Range("K136").Value = 6435.12
$ws.Range("L136").Value = 8914.667099999999
$ws.Range("M136").Value = -3885.12
$ws.Range("N136").Value = -14014.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 6003333.5
$ws.Range("I7").Value = 6003333.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 6003333.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -6003220.5
$ws.Range("N7").ClearContents()
$ws.Range("H15").Value = 31666.5
$ws.Range("J15").Value = 31666.5
$ws.Range("L15").Value = 31666.5
$ws.Range("N15").Value = -32120.5
$ws.Range("H22").Value = 373.33334
$ws.Range("I22").Value = 373.33334
$ws.Range("K22").Value = 373.33334
$ws.Range("M22").Value = -200.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 50000
$ws.Range("I12").Value = 50000
$ws.Range("K12").Value = 50000
$ws.Range("M12").Value = -49830
$ws.Range("H31").Value = 3106.7627
$ws.Range("I31").Value = 2757.5417
$ws.Range("J31").Value = 3346.2285
$ws.Range("K31").Value = 2757.5417
$ws.Range("L31").Value = 3346.2285
$ws.Range("M31").Value = -2462.5417
$ws.Range("N31").Value = -3936.2285
$ws.Range("H34").Value = 3106.7627
$ws.Range("I34").Value = 2757.5417
$ws.Range("J34").Value = 3346.2285
$ws.Range("K34").Value = 2757.5417
$ws.Range("L34").Value = 3346.2285
$ws.Range("M34").Value = -2555.5417
$ws.Range("N34").Value = -3750.2285
$ws.Range("H132").Value = 4607.75
$ws.Range("I132").Value = 4978.769
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 14936.307
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -12406.307
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 5779.6665
$ws.Range("I134").Value = 5990.5557
$ws.Range("J134").Value = 5463.3335
$ws.Range("K134").Value = 17971.6671
$ws.Range("L134").Value = 16390.0005
$ws.Range("M134").Value = -15436.6671
$ws.Range("N134").Value = -21460.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 4761947.5
$ws.Range("I12").Value = 10000050
$ws.Range("J12").Value = 36.18182
$ws.Range("K12").Value = 30000150
$ws.Range("L12").Value = 108.54546
$ws.Range("M12").Value = -29999977
$ws.Range("N12").Value = -454.54546
$ws.Range("H101").Value = 9009.666999999999
$ws.Range("J101").Value = 9009.666999999999
$ws.Range("L101").Value = 27029.001
$ws.Range("N101").Value = -31897.001
$ws.Range("H131").Value = 22001088
$ws.Range("I131").Value = 5882939.5
$ws.Range("J131").Value = 30304376
$ws.Range("K131").Value = 17648818.5
$ws.Range("L131").Value = 90913128
$ws.Range("M131").Value = -17643778.5
$ws.Range("N131").Value = -90923208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5681.518
$ws.Range("I70").Value = 5754.884
$ws.Range("J70").Value = 5438.846
$ws.Range("K70").Value = 5754.884
$ws.Range("L70").Value = 5438.846
$ws.Range("M70").Value = -5484.884
$ws.Range("N70").Value = -5978.846
$ws.Range("H73").Value = 5681.518
$ws.Range("I73").Value = 5754.884
$ws.Range("J73").Value = 5438.846
$ws.Range("K73").Value = 5754.884
$ws.Range("L73").Value = 5438.846
$ws.Range("M73").Value = -4818.884
$ws.Range("N73").Value = -7310.846
$ws.Range("H132").Value = 3540.0256
$ws.Range("I132").Value = 4388.8
$ws.Range("J132").Value = 3009.5417
$ws.Range("K132").Value = 13166.4
$ws.Range("L132").Value = 9028.625100000001
$ws.Range("M132").Value = -10636.4
$ws.Range("N132").Value = -14088.6251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 2950
$ws.Range("I107").Value = 2950
$ws.Range("K107").Value = 2950
$ws.Range("M107").Value = -1030
$ws.Range("H132").Value = 17338138
$ws.Range("I132").Value = 20639856
$ws.Range("J132").Value = 4120
$ws.Range("K132").Value = 61919568
$ws.Range("L132").Value = 12360
$ws.Range("N132").Value = -17420


Write-Host "Updated Leve profit market data across sheets."
